$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# For each year block of 4 rows (A,B,C,D quarters), swap the contents (columns A:E) of
# the 2nd and 3rd rows within the block (the "B" and "C" quarter rows), so the "C" quarter
# row moves up one row and the "B" quarter row moves down one row.
$cols = @("A","B","C","D","E")
for ($base = 2; $base -le 62; $base += 4) {
    $r1 = $base + 1
    $r2 = $base + 2
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        if ($v1 -ne $v2) {
            $ws.Range($addr1).Value = $v2
            $ws.Range($addr2).Value = $v1
        }
    }
}

# Remove the now-unneeded "产销率" (F) and "销售量" (G) columns entirely.
$ws.Range("F1:G65").EntireColumn.Delete()
